$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 592.5
$ws.Range("I6").Value2 = 592.5
$ws.Range("K6").Value2 = 1777.5
$ws.Range("M6").Value2 = -1665.5

$ws.Range("H32").Value2 = 14288741
$ws.Range("I32").Value2 = 2000
$ws.Range("J32").Value2 = 16669865
$ws.Range("K32").Value2 = 2000
$ws.Range("L32").Value2 = 16669865
$ws.Range("M32").Value2 = -1674
$ws.Range("N32").Value2 = -16670517

$ws.Range("H40").Value2 = 4027.5
$ws.Range("I40").Value2 = 4058.5454
$ws.Range("J40").Value2 = 3913.6667
$ws.Range("K40").Value2 = 4058.5454
$ws.Range("L40").Value2 = 3913.6667
$ws.Range("M40").Value2 = -3883.5454
$ws.Range("N40").Value2 = -4263.6667

$ws.Range("H80").Value2 = 2634.111
$ws.Range("J80").Value2 = 3579.389
$ws.Range("L80").Value2 = 10738.167
$ws.Range("N80").Value2 = -12734.167

$ws.Range("H83").Value2 = 2634.111
$ws.Range("J83").Value2 = 3579.389
$ws.Range("L83").Value2 = 32214.501
$ws.Range("N83").Value2 = -42198.501

$ws.Range("H113").Value2 = 10946.154
$ws.Range("I113").Value2 = 9480
$ws.Range("J113").Value2 = 11862.5
$ws.Range("K113").Value2 = 9480
$ws.Range("L113").Value2 = 11862.5
$ws.Range("M113").Value2 = -6226
$ws.Range("N113").Value2 = -18370.5

$ws.Range("H130").Value2 = 19997.143
$ws.Range("J130").Value2 = 19997.143
$ws.Range("L130").Value2 = 19997.143
$ws.Range("N130").Value2 = -30037.143

$ws.Range("H137").Value2 = 11129882
$ws.Range("I137").Value2 = 14307598
$ws.Range("J137").Value2 = 7873.75
$ws.Range("K137").Value2 = 42922794
$ws.Range("L137").Value2 = 23621.25
$ws.Range("M137").Value2 = -42920244
$ws.Range("N137").Value2 = -28721.25

$ws.Range("H138").Value2 = 4398.45
$ws.Range("J138").Value2 = 4998.231
$ws.Range("L138").Value2 = 14994.693
$ws.Range("N138").Value2 = -25274.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value2 = 3444
$ws.Range("I26").Value2 = 3444
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 3444
$ws.Range("L26").Value2 = 0
$ws.Range("M26").Value2 = -3114
$ws.Range("N26").ClearContents()

$ws.Range("H102").Value2 = 1666.6428
$ws.Range("I102").Value2 = 1487.1538
$ws.Range("J102").Value2 = 4000
$ws.Range("K102").Value2 = 1487.1538
$ws.Range("L102").Value2 = 4000
$ws.Range("M102").Value2 = 134.8462
$ws.Range("N102").Value2 = -7244

$ws.Range("H123").Value2 = 0
$ws.Range("J123").Value2 = 0
$ws.Range("L123").Value2 = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value2 = 1566
$ws.Range("I132").Value2 = 1299
$ws.Range("J132").Value2 = 2100
$ws.Range("K132").Value2 = 3897
$ws.Range("L132").Value2 = 6300
$ws.Range("M132").Value2 = -1367
$ws.Range("N132").Value2 = -11360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value2 = 1766.6666
$ws.Range("I8").Value2 = 1766.6666
$ws.Range("J8").Value2 = 0
$ws.Range("K8").Value2 = 1766.6666
$ws.Range("L8").Value2 = 0
$ws.Range("M8").Value2 = -1626.6666
$ws.Range("N8").ClearContents()

$ws.Range("H20").Value2 = 85046.914
$ws.Range("I20").Value2 = 1756.3
$ws.Range("K20").Value2 = 1756.3
$ws.Range("M20").Value2 = -1509.3

$ws.Range("H86").Value2 = 9025.692
$ws.Range("I86").Value2 = 8199.625
$ws.Range("J86").Value2 = 10347.4
$ws.Range("K86").Value2 = 8199.625
$ws.Range("L86").Value2 = 10347.4
$ws.Range("M86").Value2 = -7076.625
$ws.Range("N86").Value2 = -12593.4

$ws.Range("H89").Value2 = 9025.692
$ws.Range("I89").Value2 = 8199.625
$ws.Range("J89").Value2 = 10347.4
$ws.Range("K89").Value2 = 40998.125
$ws.Range("L89").Value2 = 51737
$ws.Range("M89").Value2 = -35382.125
$ws.Range("N89").Value2 = -62969

$ws.Range("H105").Value2 = 3255.1667
$ws.Range("I105").Value2 = 2885
$ws.Range("J105").Value2 = 3995.5
$ws.Range("K105").Value2 = 2885
$ws.Range("L105").Value2 = 3995.5
$ws.Range("M105").Value2 = -1138
$ws.Range("N105").Value2 = -7489.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 7750.1177
$ws.Range("I31").Value2 = 10093.546
$ws.Range("K31").Value2 = 10093.546
$ws.Range("M31").Value2 = -9798.546

$ws.Range("H34").Value2 = 7750.1177
$ws.Range("I34").Value2 = 10093.546
$ws.Range("K34").Value2 = 10093.546
$ws.Range("M34").Value2 = -9891.546

$ws.Range("H122").Value2 = 211664.48
$ws.Range("I122").Value2 = 271953.06
$ws.Range("J122").Value2 = 3394.9092
$ws.Range("K122").Value2 = 815859.1799999999
$ws.Range("L122").Value2 = 10184.7276
$ws.Range("M122").Value2 = -813409.1799999999
$ws.Range("N122").Value2 = -15084.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value2 = 186.68182
$ws.Range("I10").Value2 = 47.666668
$ws.Range("J10").Value2 = 812.25
$ws.Range("K10").Value2 = 143.000004
$ws.Range("L10").Value2 = 2436.75
$ws.Range("M10").Value2 = -4.00000399999999
$ws.Range("N10").Value2 = -2714.75

$ws.Range("H41").Value2 = 83333540
$ws.Range("I41").Value2 = 199
$ws.Range("K41").Value2 = 597
$ws.Range("M41").Value2 = -259

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value2 = 0
$ws.Range("I12").Value2 = 0
$ws.Range("K12").Value2 = 0
$ws.Range("M12").ClearContents()

$ws.Range("H97").Value2 = 7145.4443
$ws.Range("I97").Value2 = 1601.125
$ws.Range("K97").Value2 = 1601.125
$ws.Range("M97").Value2 = -1105.125

$ws.Range("H122").Value2 = 4694.846
$ws.Range("I122").Value2 = 4272.6875
$ws.Range("J122").Value2 = 5370.3
$ws.Range("K122").Value2 = 12818.0625
$ws.Range("L122").Value2 = 16110.9
$ws.Range("M122").Value2 = -10368.0625
$ws.Range("N122").Value2 = -21010.9

$ws.Range("H132").Value2 = 5956.4707
$ws.Range("I132").Value2 = 5930.077
$ws.Range("K132").Value2 = 17790.231
$ws.Range("M132").Value2 = -15260.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4641.478
$ws.Range("I7").Value2 = 4416.905
$ws.Range("K7").Value2 = 4416.905
$ws.Range("M7").Value2 = -4304.905

$ws.Range("H122").Value2 = 2916.6667
$ws.Range("I122").Value2 = 2500
$ws.Range("K122").Value2 = 7500
$ws.Range("M122").Value2 = -5050

$ws.Range("H126").Value2 = 4641.478
$ws.Range("I126").Value2 = 4416.905
$ws.Range("K126").Value2 = 13250.715
$ws.Range("M126").Value2 = -10780.715

$ws.Range("H132").Value2 = 3752.4167
$ws.Range("I132").Value2 = 3802.4856
$ws.Range("K132").Value2 = 11407.4568
$ws.Range("M132").Value2 = -8877.4568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 1014.2857
$ws.Range("I107").Value2 = 1049.8334
$ws.Range("J107").Value2 = 801
$ws.Range("K107").Value2 = 3149.5002
$ws.Range("L107").Value2 = 2403
$ws.Range("M107").Value2 = -1229.5002
$ws.Range("N107").Value2 = -6243

$ws.Range("H122").Value2 = 4994.75
$ws.Range("I122").Value2 = 4989.5
$ws.Range("K122").Value2 = 14968.5
$ws.Range("M122").Value2 = -12518.5

$ws.Range("H138").Value2 = 75000
$ws.Range("J138").Value2 = 75000
$ws.Range("L138").Value2 = 75000
$ws.Range("N138").Value2 = -85280

